# Update department column for the "courses" worksheet.
# Rows 2-4 (Automotive Electrical / Light Vehicle / AUR30320 courses) move
# from the old "FACULTY OF TECH SCIENCES" department label to "AUTOMOTIVE".
# Rows 5-6 (package offers) move to "Packages".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

$ws.Range("C2").Value = "AUTOMOTIVE"
$ws.Range("C3").Value = "AUTOMOTIVE"
$ws.Range("C4").Value = "AUTOMOTIVE"
$ws.Range("C5").Value = "Packages"
$ws.Range("C6").Value = "Packages"
